# Improve column matching and scoring-rule logic for the processed scorecard.
#
#  - Column matching: when looking up "the" target-value column and
#    "the" scoring-rule column for a row, prefer the *annual* variant
#    (年度目标值 / 计分规则) over the *semi-annual* one (半年度目标值 /
#    半年度计分规则) by explicitly excluding any header containing
#    "半年度" from the match.
#  - Baseline extraction (推导底线值, col O): instead of a blind ratio of
#    the semi-annual target, parse the prose scoring-rule text (计分规则,
#    col L) for its actual pass/fail score threshold (usually "低于60分，
#    不得分", sometimes 80, etc.) and its per-unit deduction, and solve for
#    the indicator value that lands exactly on that threshold. This also
#    adds ratio-based rule detection ("实际达成率/达成率目标值*100") as a
#    distinct pattern from the flat per-percent/per-unit deduction rules,
#    and recognises "每高X%，扣Y分" as an *inverse* ("lower is better")
#    rule, flipping 指标方向 (col R) to "逆向" and mirroring the generated
#    rule template (规范版计分规则, col P) accordingly.
#  - Rows whose rule text now parses successfully have 解析状态 (col Q)
#    flipped from "人工校验" (manual check needed) to "成功" (succeeded).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 1
$firstDataRow = 2
$lastCol = 18

# ---------------------------------------------------------------------
# 1. Column matching - locate columns by header text, excluding any
#    "半年度" (semi-annual) variant, so the annual columns win.
# ---------------------------------------------------------------------
$targetCol = $null      # 年度目标值
$ruleCol = $null        # 计分规则
$baselineCol = $null    # 推导底线值
$normRuleCol = $null    # 规范版计分规则
$statusCol = $null      # 解析状态
$directionCol = $null   # 指标方向

for ($c = 1; $c -le $lastCol; $c++) {
    $h = $ws.Cells.Item($headerRow, $c).Text
    if ($h -like "*半年度*") { continue }
    if ($h -eq "年度目标值" -or $h -eq "目标值")   { $targetCol    = $c }
    if ($h -eq "计分规则")                          { $ruleCol      = $c }
    if ($h -eq "推导底线值")                        { $baselineCol  = $c }
    if ($h -eq "规范版计分规则")                    { $normRuleCol  = $c }
    if ($h -eq "解析状态")                          { $statusCol    = $c }
    if ($h -eq "指标方向")                          { $directionCol = $c }
}

$lastRow = $ws.UsedRange.Rows.Count

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Format-Num([double]$n) {
    # Render doubles the way the source sheet does: plain integers with no
    # trailing decimal point ("6", not "6.0"), trimmed decimals otherwise
    # ("0.65", not "0.6500000000000001").
    $rounded = [Math]::Round($n, 10)
    if ($rounded -eq [Math]::Floor($rounded)) {
        return [string]([int64]$rounded)
    }
    return $rounded.ToString("0.##########")
}

function Get-BaselineInfo([double]$target, [string]$ruleText) {
    # Parse the free-text scoring rule (计分规则) and derive the indicator
    # value that exactly hits the rule's base/pass score threshold (底线值).
    # Returns @{ Baseline = <double>; Inverse = <bool> }, or $null if the
    # rule text doesn't match any recognised pattern.

    $threshold = 60.0
    if ($ruleText -match '低于(\d+(\.\d+)?)分') {
        $threshold = [double]$matches[1]
    }

    # Ratio-based rule: "实际达成率/达成率目标值*100" -> baseline is the
    # fraction of target that scores exactly the threshold.
    if ($ruleText -match '实际达成率.*达成率目标值.*\*\s*100') {
        return @{ Baseline = ($target * $threshold / 100.0); Inverse = $false }
    }

    # "每低X%，扣Y分" - forward indicator, percent-based deduction below target.
    if ($ruleText -match '每低(\d+(\.\d+)?)\s*%，?\s*扣\s*(\d+(\.\d+)?)\s*分') {
        $x = [double]$matches[1]
        $y = [double]$matches[3]
        $pctDrop = (100.0 - $threshold) / $y * $x
        return @{ Baseline = ($target - $pctDrop / 100.0); Inverse = $false }
    }

    # "每高X%，扣Y分" - inverse indicator ("the lower, the better"):
    # going *above* target costs points, so the baseline sits above target.
    if ($ruleText -match '每高(\d+(\.\d+)?)\s*%，?\s*扣\s*(\d+(\.\d+)?)\s*分') {
        $x = [double]$matches[1]
        $y = [double]$matches[3]
        $pctRise = (100.0 - $threshold) / $y * $x
        return @{ Baseline = ($target + $pctRise / 100.0); Inverse = $true }
    }

    # "每少X个扣Y分" - forward indicator, unit-based (not percent) deduction.
    if ($ruleText -match '每少(\d+(\.\d+)?)\s*个扣\s*(\d+(\.\d+)?)\s*分') {
        $x = [double]$matches[1]
        $y = [double]$matches[3]
        $unitDrop = (100.0 - $threshold) / $y * $x
        return @{ Baseline = ($target - $unitDrop); Inverse = $false }
    }

    return $null
}

function New-NormalRuleText([string]$target, [string]$baseline) {
    return "P为指标实际值，$target" + "为目标值，$baseline" + "为底线值。`n" + `
        "1.若P≥$target" + "，得100分（满分）；`n" + `
        "2.若$baseline" + "<P<$target" + "，按线性比例计算，即：得分=60+(P-$baseline)/($target-$baseline)×(100-60)；`n" + `
        "3.若P=$baseline" + "，得60分（基础分）；`n" + `
        "4.若P<$baseline" + "，得0分。"
}

function New-InverseRuleText([string]$target, [string]$baseline) {
    return "P为指标实际值，$target" + "为目标值，$baseline" + "为底线值。`n" + `
        "1.若P≤$target" + "，得100分（满分）；`n" + `
        "2.若$target" + "<P<$baseline" + "，按线性比例计算，即：得分=100-(P-$target)/($baseline-$target)×(100-60)；`n" + `
        "3.若P=$baseline" + "，得60分（基础分）；`n" + `
        "4.若P＞$baseline" + "，得0分。"
}

function Set-TextValue($range, [string]$text) {
    # Force the cell to stay a TEXT cell (matches the rest of column O,
    # which stores numeric-looking baselines as text) instead of letting
    # Excel auto-convert a numeric-looking string into a real Number. The
    # leading apostrophe is Excel's normal "treat as text" quote-prefix;
    # resetting the style back to Normal afterwards clears the cosmetic
    # quote-prefix marker so the cell's style is unaffected.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2. Walk every data row, recompute the baseline from the annual target +
#    scoring-rule columns, and only touch rows whose parse now succeeds.
# ---------------------------------------------------------------------
for ($row = $firstDataRow; $row -le $lastRow; $row++) {
    $targetVal = $ws.Cells.Item($row, $targetCol).Value2
    $ruleText = $ws.Cells.Item($row, $ruleCol).Text

    if ([string]::IsNullOrEmpty($ruleText)) { continue }

    $isNumericTarget = ($targetVal -is [double]) -or ($targetVal -is [int])
    $target = 0.0
    if ($isNumericTarget) { $target = [double]$targetVal }

    $info = Get-BaselineInfo $target $ruleText
    if ($info -eq $null) { continue }

    $targetText = Format-Num $target
    $baselineText = Format-Num $info.Baseline

    # Column O: 推导底线值
    Set-TextValue $ws.Cells.Item($row, $baselineCol) $baselineText

    # Column P: 规范版计分规则
    if ($info.Inverse) {
        $ruleOut = New-InverseRuleText $targetText $baselineText
    } else {
        $ruleOut = New-NormalRuleText $targetText $baselineText
    }
    $ws.Cells.Item($row, $normRuleCol).Value = $ruleOut

    # Column Q: 解析状态 - parsing succeeded.
    $ws.Cells.Item($row, $statusCol).Value = "成功"

    # Column R: 指标方向 - flip to inverse only when the rule calls for it.
    if ($info.Inverse) {
        $ws.Cells.Item($row, $directionCol).Value = "逆向"
    }
}

Write-Host "Recomputed baseline/rule columns using annual-target column matching"
